$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A19").Value = "t"
